$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Cells.Item(2, 3).Value = 7.719167780743623
$ws.Cells.Item(2, 4).Value = 6.733769805988533
$ws.Cells.Item(2, 5).Value = 9.252933337777126
$ws.Cells.Item(2, 6).Value = 61.39123680782699
$ws.Cells.Item(2, 7).Value = 3.757078435056697
$ws.Cells.Item(2, 9).Value = 45.26254811337704
$ws.Cells.Item(2, 10).Value = 8.661638453971412
$ws.Cells.Item(2, 13).Value = 35.84339797314053
$ws.Cells.Item(2, 14).Value = 17.02528174367247
$ws.Cells.Item(3, 3).Value = 7.73116687017151
$ws.Cells.Item(3, 4).Value = 6.407882797737234
$ws.Cells.Item(3, 5).Value = 9.04764189062735
$ws.Cells.Item(3, 6).Value = 60.75512048053122
$ws.Cells.Item(3, 7).Value = 3.766504137164751
$ws.Cells.Item(3, 9).Value = 44.74718523391278
$ws.Cells.Item(3, 10).Value = 8.689269159263574
$ws.Cells.Item(3, 13).Value = 35.00617729778023
$ws.Cells.Item(3, 14).Value = 16.90437494280274
$ws.Cells.Item(4, 3).Value = 7.74141508071552
$ws.Cells.Item(4, 4).Value = 6.202125597637469
$ws.Cells.Item(4, 5).Value = 8.918398624291571
$ws.Cells.Item(4, 6).Value = 60.38964547714774
$ws.Cells.Item(4, 7).Value = 3.772557799243368
$ws.Cells.Item(4, 9).Value = 44.44855158456991
$ws.Cells.Item(4, 10).Value = 8.70710312606283
$ws.Cells.Item(4, 13).Value = 34.49017388580899
$ws.Cells.Item(4, 14).Value = 16.83224193979358
$ws.Cells.Item(5, 3).Value = 7.74630601846193
$ws.Cells.Item(5, 4).Value = 6.116978202065504
$ws.Cells.Item(5, 5).Value = 8.864965892572387
$ws.Cells.Item(5, 6).Value = 60.24710793984875
$ws.Cells.Item(5, 7).Value = 3.775092173187493
$ws.Cells.Item(5, 9).Value = 44.33141546828637
$ws.Cells.Item(5, 10).Value = 8.714589999584177
$ws.Cells.Item(5, 13).Value = 34.27969418581768
$ws.Cells.Item(5, 14).Value = 16.80340323175375
$ws.Cells.Item(6, 3).Value = 7.747161034651355
$ws.Cells.Item(6, 4).Value = 6.102765020432185
$ws.Cells.Item(6, 5).Value = 8.856048373078275
$ws.Cells.Item(6, 6).Value = 60.22382815162215
$ws.Cells.Item(6, 7).Value = 3.775517093095182
$ws.Cells.Item(6, 9).Value = 44.31224267923453
$ws.Cells.Item(6, 10).Value = 8.715846467960924
$ws.Cells.Item(6, 13).Value = 34.24474038299393
$ws.Cells.Item(6, 14).Value = 16.79864894600812
$ws.Cells.Item(7, 3).Value = 7.741478160704049
$ws.Cells.Item(7, 4).Value = 6.20098235542252
$ws.Cells.Item(7, 5).Value = 8.917681057702998
$ws.Cells.Item(7, 6).Value = 60.38769716800528
$ws.Cells.Item(7, 7).Value = 3.772591704892084
$ws.Cells.Item(7, 9).Value = 44.44695328363666
$ws.Cells.Item(7, 10).Value = 8.707203207079131
$ws.Cells.Item(7, 13).Value = 34.48733572553061
$ws.Cells.Item(7, 14).Value = 16.83185072424887
$ws.Cells.Item(8, 3).Value = 7.722701419953141
$ws.Cells.Item(8, 4).Value = 6.622647828723341
$ws.Cells.Item(8, 5).Value = 9.182833746076811
$ws.Cells.Item(8, 6).Value = 61.16674925497711
$ws.Cells.Item(8, 7).Value = 3.760273467120116
$ws.Cells.Item(8, 9).Value = 45.08120995575432
$ws.Cells.Item(8, 10).Value = 8.670985831236507
$ws.Cells.Item(8, 13).Value = 35.555296406234
$ws.Cells.Item(8, 14).Value = 16.98317269628054
$ws.Cells.Item(9, 3).Value = 7.709160327638885
$ws.Cells.Item(9, 4).Value = 7.529804234442093
$ws.Cells.Item(9, 5).Value = 9.675909704109579
$ws.Cells.Item(9, 6).Value = 62.88919142733415
$ws.Cells.Item(9, 7).Value = 3.738206520038428
$ws.Cells.Item(9, 9).Value = 46.46201891434184
$ws.Cells.Item(9, 10).Value = 8.606810879874899
$ws.Cells.Item(9, 13).Value = 37.62213675707093
$ws.Cells.Item(9, 14).Value = 17.29540048636253
$ws.Cells.Item(10, 3).Value = 7.713967175810203
$ws.Cells.Item(10, 4).Value = 8.166409614226787
$ws.Cells.Item(10, 5).Value = 10.01988239463606
$ws.Cells.Item(10, 6).Value = 64.26646572996549
$ws.Cells.Item(10, 7).Value = 3.72323428119177
$ws.Cells.Item(10, 9).Value = 47.55373666159723
$ws.Cells.Item(10, 10).Value = 8.563772363081849
$ws.Cells.Item(10, 13).Value = 39.10881434249271
$ws.Cells.Item(10, 14).Value = 17.53258297322792
$ws.Cells.Item(11, 3).Value = 7.719472902337553
$ws.Cells.Item(11, 4).Value = 8.441133479620184
$ws.Cells.Item(11, 5).Value = 10.1720403671635
$ws.Cells.Item(11, 6).Value = 64.91559114158453
$ws.Cells.Item(11, 7).Value = 3.716685025177187
$ws.Cells.Item(11, 9).Value = 48.06566139902349
$ws.Cells.Item(11, 10).Value = 8.545071789020058
$ws.Cells.Item(11, 13).Value = 39.77533544944738
$ws.Cells.Item(11, 14).Value = 17.64181705169979
$ws.Cells.Item(12, 3).Value = 7.722043960760952
$ws.Cells.Item(12, 4).Value = 8.543036762625757
$ws.Cells.Item(12, 5).Value = 10.22901150706671
$ws.Cells.Item(12, 6).Value = 65.16449099293791
$ws.Cells.Item(12, 7).Value = 3.714242026924496
$ws.Cells.Item(12, 9).Value = 48.26158021613972
$ws.Cells.Item(12, 10).Value = 8.538115503364535
$ws.Cells.Item(12, 13).Value = 40.02610465953984
$ws.Cells.Item(12, 14).Value = 17.68334337503706
$ws.Cells.Item(13, 3).Value = 7.721468477418732
$ws.Cells.Item(13, 4).Value = 8.521184567492318
$ws.Cells.Item(13, 5).Value = 10.21677096705595
$ws.Cells.Item(13, 6).Value = 65.11075111417705
$ws.Cells.Item(13, 7).Value = 3.714766531098482
$ws.Cells.Item(13, 9).Value = 48.2192959826928
$ws.Cells.Item(13, 10).Value = 8.539608111327649
$ws.Cells.Item(13, 13).Value = 39.97217277465219
$ws.Cells.Item(13, 14).Value = 17.67439323303813
$ws.Cells.Item(14, 3).Value = 7.719674631823663
$ws.Cells.Item(14, 4).Value = 8.449559723894946
$ws.Cells.Item(14, 5).Value = 10.17674054559519
$ws.Cells.Item(14, 6).Value = 64.93600728268734
$ws.Cells.Item(14, 7).Value = 3.716483298744047
$ws.Cells.Item(14, 9).Value = 48.0817391783466
$ws.Cells.Item(14, 10).Value = 8.54449698760215
$ws.Cells.Item(14, 13).Value = 39.79600021861461
$ws.Cells.Item(14, 14).Value = 17.64523038143216
$ws.Cells.Item(15, 3).Value = 7.718639419774312
$ws.Cells.Item(15, 4).Value = 8.405410541057549
$ws.Cells.Item(15, 5).Value = 10.15213563652355
$ws.Cells.Item(15, 6).Value = 64.82936924764398
$ws.Cells.Item(15, 7).Value = 3.717539678114462
$ws.Cells.Item(15, 9).Value = 47.99774638862268
$ws.Cells.Item(15, 10).Value = 8.547507841764874
$ws.Cells.Item(15, 13).Value = 39.68787114953473
$ws.Cells.Item(15, 14).Value = 17.62738744617089
$ws.Cells.Item(16, 3).Value = 7.713674918549815
$ws.Cells.Item(16, 4).Value = 8.148157313462422
$ws.Cells.Item(16, 5).Value = 10.00984928039018
$ws.Cells.Item(16, 6).Value = 64.22448508269815
$ws.Cells.Item(16, 7).Value = 3.723667508076232
$ws.Cells.Item(16, 9).Value = 47.52057746144705
$ws.Cells.Item(16, 10).Value = 8.565012067575957
$ws.Cells.Item(16, 13).Value = 39.06503952230999
$ws.Cells.Item(16, 14).Value = 17.52546857614799
$ws.Cells.Item(17, 3).Value = 7.711485874661673
$ws.Cells.Item(17, 4).Value = 7.986538323729131
$ws.Cells.Item(17, 5).Value = 9.921435539764783
$ws.Cells.Item(17, 6).Value = 63.85908442393978
$ws.Cells.Item(17, 7).Value = 3.727493358070421
$ws.Cells.Item(17, 9).Value = 47.23167218032239
$ws.Cells.Item(17, 10).Value = 8.575974463772743
$ws.Cells.Item(17, 13).Value = 38.6802851237337
$ws.Cells.Item(17, 14).Value = 17.46326561006358
$ws.Cells.Item(18, 3).Value = 7.710538690347814
$ws.Cells.Item(18, 4).Value = 7.892178908713403
$ws.Cells.Item(18, 5).Value = 9.870177422810539
$ws.Cells.Item(18, 6).Value = 63.65105249838352
$ws.Cells.Item(18, 7).Value = 3.729718552006556
$ws.Cells.Item(18, 9).Value = 47.06694969029641
$ws.Cells.Item(18, 10).Value = 8.582362434580162
$ws.Cells.Item(18, 13).Value = 38.45807660200021
$ws.Cells.Item(18, 14).Value = 17.42761594216913
$ws.Cells.Item(19, 3).Value = 7.710271282364772
$ws.Cells.Item(19, 4).Value = 7.859989635469742
$ws.Cells.Item(19, 5).Value = 9.852753631394551
$ws.Cells.Item(19, 6).Value = 63.58098822318007
$ws.Cells.Item(19, 7).Value = 3.730476217997278
$ws.Cells.Item(19, 9).Value = 47.01143029487923
$ws.Cells.Item(19, 10).Value = 8.584539524739426
$ws.Cells.Item(19, 13).Value = 38.38269172085024
$ws.Cells.Item(19, 14).Value = 17.41556846989702
$ws.Cells.Item(20, 3).Value = 7.711686556701074
$ws.Cells.Item(20, 4).Value = 8.003887816747216
$ws.Cells.Item(20, 5).Value = 9.930889439548183
$ws.Cells.Item(20, 6).Value = 63.89776192471783
$ws.Cells.Item(20, 7).Value = 3.727083541025959
$ws.Cells.Item(20, 9).Value = 47.26227773222584
$ws.Cells.Item(20, 10).Value = 8.574798947179001
$ws.Cells.Item(20, 13).Value = 38.72133852139552
$ws.Cells.Item(20, 14).Value = 17.46987421460592
$ws.Cells.Item(21, 3).Value = 7.720188264415562
$ws.Cells.Item(21, 4).Value = 8.470655362299361
$ws.Cells.Item(21, 5).Value = 10.18851623450745
$ws.Cells.Item(21, 6).Value = 64.9872512121387
$ws.Cells.Item(21, 7).Value = 3.715978041177951
$ws.Cells.Item(21, 9).Value = 48.12208802322039
$ws.Cells.Item(21, 10).Value = 8.543057615999636
$ws.Cells.Item(21, 13).Value = 39.84779223248614
$ws.Cells.Item(21, 14).Value = 17.65379207261597
$ws.Cells.Item(22, 3).Value = 7.728582211324252
$ws.Cells.Item(22, 4).Value = 8.763313026046434
$ws.Cells.Item(22, 5).Value = 10.35310366098413
$ws.Cells.Item(22, 6).Value = 65.71722990792533
$ws.Cells.Item(22, 7).Value = 3.708935713638599
$ws.Cells.Item(22, 9).Value = 48.69599006405193
$ws.Cells.Item(22, 10).Value = 8.523042191586896
$ws.Cells.Item(22, 13).Value = 40.57441566408314
$ws.Cells.Item(22, 14).Value = 17.77492164221672
$ws.Cells.Item(23, 3).Value = 7.723839717157949
$ws.Cells.Item(23, 4).Value = 8.608247243178429
$ws.Cells.Item(23, 5).Value = 10.26561501816244
$ws.Cells.Item(23, 6).Value = 65.32603944993664
$ws.Cells.Item(23, 7).Value = 3.712674785206808
$ws.Cells.Item(23, 9).Value = 48.3886382449435
$ws.Cells.Item(23, 10).Value = 8.533658403522841
$ws.Cells.Item(23, 13).Value = 40.18754843694864
$ws.Cells.Item(23, 14).Value = 17.71019767115853
$ws.Cells.Item(24, 3).Value = 7.711594859316108
$ws.Cells.Item(24, 4).Value = 7.996048604755752
$ws.Cells.Item(24, 5).Value = 9.926616662129634
$ws.Cells.Item(24, 6).Value = 63.88026945490775
$ws.Cells.Item(24, 7).Value = 3.727268739337976
$ws.Cells.Item(24, 9).Value = 47.24843667534476
$ws.Cells.Item(24, 10).Value = 8.575330131706076
$ws.Cells.Item(24, 13).Value = 38.70278140325364
$ws.Cells.Item(24, 14).Value = 17.46688611440916
$ws.Cells.Item(25, 3).Value = 7.710272028456379
$ws.Cells.Item(25, 4).Value = 7.282127455571493
$ws.Cells.Item(25, 5).Value = 9.545603500828873
$ws.Cells.Item(25, 6).Value = 62.40304315440495
$ws.Cells.Item(25, 7).Value = 3.743955962449536
$ws.Cells.Item(25, 9).Value = 46.07442038084357
$ws.Cells.Item(25, 10).Value = 8.623445378326547
$ws.Cells.Item(25, 13).Value = 37.06755431529721
$ws.Cells.Item(25, 14).Value = 17.20944582058611
